# Subject 17-20: populate remaining response rows (14-32) and add a new
# "Valid" column (R) marking whether each subject's responses are complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header for column R
$ws.Cells.Item(1, 18).Value = "Valid"

# Existing subjects (rows 2-13) are all valid
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 18).Value = 1
}

# Full question responses (Q1..Q15) for subjects 12-20 (rows 14-22).
# Each entry: Condition, Q1..Q15, Valid
$fullRows = @{
    14 = @("G", 5,4,5,4,5,5,4,6,4,4,5,5,4,4,6, 1)
    15 = @("F", 4,3,3,4,3,2,1,3,1,2,2,3,2,3,4, 0)
    16 = @("G", 2,3,2,2,2,2,3,3,2,2,3,3,3,3,4, 0)
    17 = @("F", 5,4,4,5,4,5,3,6,2,4,4,6,4,4,4, 1)
    18 = @("G", 3,2,4,2,4,3,3,3,4,2,2,3,4,2,3, 1)
    19 = @("F", 5,4,5,2,5,4,5,4,5,2,5,2,6,4,5, 1)
    20 = @("G", 4,3,3,3,4,3,5,6,5,3,5,4,6,4,6, 1)
    21 = @("F", 1,2,3,1,3,2,4,5,6,2,6,4,6,5,5, 1)
    22 = @("G", 1,2,1,2,2,2,3,5,4,2,3,3,3,5,3, 1)
}

foreach ($r in $fullRows.Keys) {
    $data = $fullRows[$r]
    $ws.Cells.Item($r, 1).Value = $r - 2     # ID = row - 2
    $ws.Cells.Item($r, 2).Value = $data[0]   # Condition (G/F)
    for ($q = 0; $q -lt 15; $q++) {
        $ws.Cells.Item($r, 3 + $q).Value = $data[1 + $q]
    }
    $ws.Cells.Item($r, 18).Value = $data[16]
}

# Remaining subjects (21-30, rows 23-32) only have ID/Condition/Valid recorded.
$sparseRows = @{
    23 = "F"; 24 = "G"; 25 = "F"; 26 = "G"; 27 = "F"
    28 = "G"; 29 = "F"; 30 = "G"; 31 = "F"; 32 = "G"
}

foreach ($r in $sparseRows.Keys) {
    $ws.Cells.Item($r, 1).Value = $r - 2
    $ws.Cells.Item($r, 2).Value = $sparseRows[$r]
    $ws.Cells.Item($r, 18).Value = 1
}

# Trailing Valid marker with no associated subject data
$ws.Cells.Item(33, 18).Value = 1

$ws.Range("K23").Select()
